$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cluster column (D) values for the 5 data rows (row 2..6)
$targets = @("ECs", "FAPs", "M1", "M2", "sCs")

# Numeric columns E..T for each of the 5 rows, in order
$data = @(
    @(3, 1, 2.29025, 6.870749999999999, 1, 1, 1, 0.3333333333333333, 0.2365093333333333, 0.709528, 0.006013533229697845, 0.006013533229697844, 0.5416655006666666, 4.874989505999999, 0.006013533229697845, 0.006013533229697844),
    @(3, 1, 2.29025, 6.870749999999999, 1, 1, 1, 0.3333333333333333, 0.075101, 0.225303, 0.001909532925058086, 0.001909532925058085, 0.17200006525, 1.54800058725, 0.001909532925058086, 0.001909532925058085),
    @(3, 1, 2.29025, 6.870749999999999, 1, 1, 1, 0.3333333333333333, 0.006358, 0.019074, 0.0001616597693442072, 0.0001616597693442072, 0.0145614095, 0.1310526855, 0.0001616597693442072, 0.0001616597693442072),
    @(3, 1, 2.29025, 6.870749999999999, 1, 1, 2, 0.6666666666666666, 0.1007303333333333, 0.302191, 0.002561189438916605, 0.002561189438916605, 0.2306976459166666, 2.07627881325, 0.002561189438916605, 0.002561189438916605),
    @(3, 1, 2.29025, 6.870749999999999, 1, 1, 3, 1, 38.91081433333333, 116.732443, 0.9893540846369834, 0.9893540846369833, 89.11549252691667, 802.0394327422499, 0.9893540846369834, 0.9893540846369833)
)

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = "ECs"
    $ws.Cells.Item($r, 2).Value = "Fgf9"
    $ws.Cells.Item($r, 3).Value = "Fgfr4"
    $ws.Cells.Item($r, 4).Value = $targets[$i]

    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 5).Value = $row[$c]
    }
}
